# Add new daily flight rows (2023-12-01 .. 2023-12-24) to the bottom of the
# "DSM Scheduled Flights vs actual" data table, continuing the existing
# DateTime / Scheduled flights / Tracked flights / Percent columns, plus a
# trailing row that only contains the (now #DIV/0!) Percent formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data (row, date, scheduled, tracked) ---------------------------
$rows = @(
    @(1334, "2023-12-01", 73, 71),
    @(1335, "2023-12-02", 53, 52),
    @(1336, "2023-12-03", 60, 59),
    @(1337, "2023-12-04", 59, 59),
    @(1338, "2023-12-05", 65, 62),
    @(1339, "2023-12-06", 71, 70),
    @(1340, "2023-12-07", 80, 78),
    @(1341, "2023-12-08", 80, 77),
    @(1342, "2023-12-09", 60, 57),
    @(1343, "2023-12-10", 61, 61),
    @(1344, "2023-12-11", 83, 79),
    @(1345, "2023-12-12", 74, 68),
    @(1346, "2023-12-13", 79, 70),
    @(1347, "2023-12-14", 90, 87),
    @(1348, "2023-12-15", 78, 74),
    @(1349, "2023-12-16", 54, 47),
    @(1350, "2023-12-17", 64, 60),
    @(1351, "2023-12-18", 70, 67),
    @(1352, "2023-12-19", 66, 66),
    @(1353, "2023-12-20", 66, 63),
    @(1354, "2023-12-21", 93, 83),
    @(1355, "2023-12-22", 69, 66),
    @(1356, "2023-12-23", 53, 51),
    @(1357, "2023-12-24", 58, 55)
)

$firstRow = 1334
$lastDataRow = 1357
$trailingRow = 1358   # extra row with only a (div/0) Percent formula

# --- Column A: dates as literal text (matching existing "s" style) ------
# A leading apostrophe forces these to be stored as text instead of being
# auto-converted into date serials, matching the existing column A cells.
foreach ($r in $rows) {
    $rowNum = $r[0]
    $dateText = $r[1]
    $ws.Range("A$rowNum").Value = "'" + $dateText
}

# --- Columns B & C: scheduled / tracked flight counts --------------------
foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("B$rowNum").Value = $r[2]
    $ws.Range("C$rowNum").Value = $r[3]
}

# --- Column D: Percent = Tracked / Scheduled -----------------------------
# Written as two contiguous formula blocks (mirrors the source workbook's
# shared-formula grouping), the first covering the new data rows and the
# second covering the data rows plus the trailing, data-less row whose
# B/C are blank and therefore evaluates to #DIV/0!.
$ws.Range("D$firstRow`:D1345").Formula = "=C$firstRow/B$firstRow"
$ws.Range("D1346:D$trailingRow").Formula = "=C1346/B1346"

# --- Copy number formats / styles from the last pre-existing data row ---
# so the newly written cells pick up the same direct formatting (s indexes)
# as the rest of the table, rather than Excel's generic defaults.
$ws.Range("A1333").Copy() | Out-Null
$ws.Range("A$firstRow`:A$lastDataRow").PasteSpecial(-4122) | Out-Null

$ws.Range("B1333:C1333").Copy() | Out-Null
$ws.Range("B$firstRow`:C$lastDataRow").PasteSpecial(-4122) | Out-Null

$ws.Range("D1333").Copy() | Out-Null
$ws.Range("D$firstRow`:D$trailingRow").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Update the view so the new bottom rows are shown / selected --------
$ws.Range("D1357:D1358").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 1336
$win.ScrollColumn = 1
